# komplex update: fooldal css fix + sajat html, css + szinek modositasa
#
# - Adds two new rows ("h1 színe" / "h2 színe") under the "Színek" section
# - Shifts the "Betűtípusok" and "Betűméretek" tables down by one row
# - Updates the h1 font size value from 30px to 55px
# - Touches F13 (an extra, empty helper cell next to the "Betűtípusok"
#   table header row) so the used range grows to column F
# - Updates the selected cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a single blank row above the "Betűtípusok" header (old row 11),
# which pushes everything from there on down by one row and turns the
# previously-empty rows 9/10 into the new "szine" rows.
$ws.Rows("11:11").Insert()

# New content for the previously blank rows 9 and 10.
$ws.Cells.Item(9, 2).Value = "h1 színe"
$ws.Cells.Item(10, 2).Value = "h2 színe"

# The "Betűméretek" table moved down by one row; h1's size changed.
$ws.Cells.Item(20, 3).Value = "55px"

# Touch the helper cell to the right of the "Betűtípusok" header row so the
# worksheet's used range extends to column F (matches the widened table).
$f13 = $ws.Cells.Item(13, 6)
$f13.VerticalAlignment = -4107
$f13.Interior.Pattern = -4142

# Update the active selection.
[void]$ws.Range("H14").Select()
